$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap row 7 and row 10 ---
$ws.Range("B7").Value = 7126860
$ws.Range("B10").Value = 7126859
$ws.Range("E7").Value = 'NK Tomislav'
$ws.Range("E10").Value = 'NK Granicar Zupanja'
$ws.Range("F7").Value = 'NK Oriolik Oriovac'
$ws.Range("F10").Value = 'Slavija Pleternica'
$ws.Range("G7").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I7").Value = 0
$ws.Range("I10").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K7").Value = 'D'
$ws.Range("K10").Value = 'D'
$ws.Range("L7").Value = 2.5
$ws.Range("L10").Value = 1.5
$ws.Range("M7").Value = 3.4
$ws.Range("M10").Value = 4
$ws.Range("N7").Value = 2.4
$ws.Range("N10").Value = 5
$ws.Range("O7").Value = 2.625
$ws.Range("O10").Value = 1.5
$ws.Range("P7").Value = 3.4
$ws.Range("P10").Value = 4
$ws.Range("Q7").Value = 2.3
$ws.Range("Q10").Value = 5
$ws.Range("R7").Value = 0
$ws.Range("R10").Value = -1
$ws.Range("S7").Value = 1.975
$ws.Range("S10").Value = 1.8
$ws.Range("T7").Value = 1.725
$ws.Range("T10").Value = 2
$ws.Range("U7").Value = 3
$ws.Range("U10").Value = 3
$ws.Range("V7").Value = 1.95
$ws.Range("V10").Value = 1.85
$ws.Range("W7").Value = 1.85
$ws.Range("W10").Value = 1.95
$ws.Range("X7").Value = -1
$ws.Range("X10").Value = -1
$ws.Range("Y7").Value = 2.4
$ws.Range("Y10").Value = 3
$ws.Range("Z7").Value = -1
$ws.Range("Z10").Value = -1
$ws.Range("AA7").Value = 0
$ws.Range("AA10").Value = -1
$ws.Range("AB7").Value = 0
$ws.Range("AB10").Value = 1
$ws.Range("AC7").Value = -1
$ws.Range("AC10").Value = -1
$ws.Range("AD7").Value = 0.8500000000000001
$ws.Range("AD10").Value = 0.95
# --- swap row 23 and row 24 ---
$ws.Range("B23").Value = 7202435
$ws.Range("B24").Value = 7202437
$ws.Range("E23").Value = 'NK Udarnik Kurilovec'
$ws.Range("E24").Value = 'NK Maksimir'
$ws.Range("F23").Value = 'NK Mladost Petrinja'
$ws.Range("F24").Value = 'Sava Strmec'
$ws.Range("G23").Value = 6
$ws.Range("G24").Value = 6
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I23").Value = 3
$ws.Range("I24").Value = 3
$ws.Range("J23").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K23").Value = 'H'
$ws.Range("K24").Value = 'H'
$ws.Range("L23").Value = 2
$ws.Range("L24").Value = 1.4
$ws.Range("M23").Value = 3.4
$ws.Range("M24").Value = 4.333
$ws.Range("N23").Value = 3.1
$ws.Range("N24").Value = 6
$ws.Range("O23").Value = 2
$ws.Range("O24").Value = 1.4
$ws.Range("P23").Value = 3.4
$ws.Range("P24").Value = 4.333
$ws.Range("Q23").Value = 3.1
$ws.Range("Q24").Value = 6
$ws.Range("R23").Value = -0.25
$ws.Range("R24").Value = -1.25
$ws.Range("S23").Value = 1.8
$ws.Range("S24").Value = 1.85
$ws.Range("T23").Value = 2
$ws.Range("T24").Value = 1.95
$ws.Range("U23").Value = 2.75
$ws.Range("U24").Value = 3
$ws.Range("V23").Value = 1.825
$ws.Range("V24").Value = 1.8
$ws.Range("W23").Value = 1.975
$ws.Range("W24").Value = 2
$ws.Range("X23").Value = 1
$ws.Range("X24").Value = 0.3999999999999999
$ws.Range("Y23").Value = -1
$ws.Range("Y24").Value = -1
$ws.Range("Z23").Value = -1
$ws.Range("Z24").Value = -1
$ws.Range("AA23").Value = 0.8
$ws.Range("AA24").Value = 0.8500000000000001
$ws.Range("AB23").Value = -1
$ws.Range("AB24").Value = -1
$ws.Range("AC23").Value = 0.825
$ws.Range("AC24").Value = 0.8
$ws.Range("AD23").Value = -1
$ws.Range("AD24").Value = -1
# --- swap row 29 and row 30 ---
$ws.Range("B29").Value = 7250138
$ws.Range("B30").Value = 7250137
$ws.Range("E29").Value = 'NK Tomislav'
$ws.Range("E30").Value = 'NK Granicar Zupanja'
$ws.Range("F29").Value = 'Sloga Nova Gradiska'
$ws.Range("F30").Value = 'NK Svacic'
$ws.Range("G29").Value = 2
$ws.Range("G30").Value = 0
$ws.Range("H29").Value = 2
$ws.Range("H30").Value = 2
$ws.Range("I29").Value = 1
$ws.Range("I30").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("J30").Value = 1
$ws.Range("K29").Value = 'D'
$ws.Range("K30").Value = 'A'
$ws.Range("L29").Value = 2.1
$ws.Range("L30").Value = 1.727
$ws.Range("M29").Value = 3.4
$ws.Range("M30").Value = 3.75
$ws.Range("N29").Value = 2.9
$ws.Range("N30").Value = 3.75
$ws.Range("O29").Value = 2.1
$ws.Range("O30").Value = 1.727
$ws.Range("P29").Value = 3.4
$ws.Range("P30").Value = 3.75
$ws.Range("Q29").Value = 2.9
$ws.Range("Q30").Value = 3.75
$ws.Range("R29").Value = -0.25
$ws.Range("R30").Value = -0.75
$ws.Range("S29").Value = 1.9
$ws.Range("S30").Value = 1.975
$ws.Range("T29").Value = 1.9
$ws.Range("T30").Value = 1.825
$ws.Range("U29").Value = 3
$ws.Range("U30").Value = 2.5
$ws.Range("V29").Value = 1.9
$ws.Range("V30").Value = 1.8
$ws.Range("W29").Value = 1.9
$ws.Range("W30").Value = 2
$ws.Range("X29").Value = -1
$ws.Range("X30").Value = -1
$ws.Range("Y29").Value = 2.4
$ws.Range("Y30").Value = -1
$ws.Range("Z29").Value = -1
$ws.Range("Z30").Value = 2.75
$ws.Range("AA29").Value = -0.5
$ws.Range("AA30").Value = -1
$ws.Range("AB29").Value = 0.45
$ws.Range("AB30").Value = 0.825
$ws.Range("AC29").Value = 0.8999999999999999
$ws.Range("AC30").Value = -1
$ws.Range("AD29").Value = -1
$ws.Range("AD30").Value = 1
# --- swap row 78 and row 79 ---
$ws.Range("B78").Value = 7519479
$ws.Range("B79").Value = 7519478
$ws.Range("E78").Value = 'Sava Strmec'
$ws.Range("E79").Value = 'NK Bistra'
$ws.Range("F78").Value = 'NK Ponikve'
$ws.Range("F79").Value = 'Lucko'
$ws.Range("G78").Value = 0
$ws.Range("G79").Value = 1
$ws.Range("H78").Value = 1
$ws.Range("H79").Value = 2
$ws.Range("I78").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J78").Value = 1
$ws.Range("J79").Value = 0
$ws.Range("K78").Value = 'A'
$ws.Range("K79").Value = 'A'
$ws.Range("L78").Value = 2.75
$ws.Range("L79").Value = 3
$ws.Range("M78").Value = 3.4
$ws.Range("M79").Value = 3.6
$ws.Range("N78").Value = 2.2
$ws.Range("N79").Value = 2
$ws.Range("O78").Value = 2.75
$ws.Range("O79").Value = 3
$ws.Range("P78").Value = 3.4
$ws.Range("P79").Value = 3.6
$ws.Range("Q78").Value = 2.2
$ws.Range("Q79").Value = 2
$ws.Range("R78").Value = 0.25
$ws.Range("R79").Value = 0.25
$ws.Range("S78").Value = 1.8
$ws.Range("S79").Value = 2
$ws.Range("T78").Value = 2
$ws.Range("T79").Value = 1.8
$ws.Range("U78").Value = 2.75
$ws.Range("U79").Value = 2.5
$ws.Range("V78").Value = 1.9
$ws.Range("V79").Value = 1.8
$ws.Range("W78").Value = 1.9
$ws.Range("W79").Value = 2
$ws.Range("X78").Value = -1
$ws.Range("X79").Value = -1
$ws.Range("Y78").Value = -1
$ws.Range("Y79").Value = -1
$ws.Range("Z78").Value = 1.2
$ws.Range("Z79").Value = 1
$ws.Range("AA78").Value = -1
$ws.Range("AA79").Value = -1
$ws.Range("AB78").Value = 1
$ws.Range("AB79").Value = 0.8
$ws.Range("AC78").Value = -1
$ws.Range("AC79").Value = 0.8
$ws.Range("AD78").Value = 0.8999999999999999
$ws.Range("AD79").Value = -1
# --- swap row 127 and row 129 ---
$ws.Range("B127").Value = 8163883
$ws.Range("B129").Value = 8163880
$ws.Range("E127").Value = 'NK Zelina'
$ws.Range("E129").Value = 'NK Dinamo Odranski Obre'
$ws.Range("F127").Value = 'NK Mladost Petrinja'
$ws.Range("F129").Value = 'Sava Strmec'
$ws.Range("G127").Value = 2
$ws.Range("G129").Value = 2
$ws.Range("H127").Value = 1
$ws.Range("H129").Value = 2
$ws.Range("I127").Value = 0
$ws.Range("I129").Value = 1
$ws.Range("J127").Value = 1
$ws.Range("J129").Value = 1
$ws.Range("K127").Value = 'H'
$ws.Range("K129").Value = 'D'
$ws.Range("L127").Value = 2
$ws.Range("L129").Value = 1.909
$ws.Range("M127").Value = 3.3
$ws.Range("M129").Value = 3.4
$ws.Range("N127").Value = 3.25
$ws.Range("N129").Value = 3.4
$ws.Range("O127").Value = 2.05
$ws.Range("O129").Value = 1.8
$ws.Range("P127").Value = 3.5
$ws.Range("P129").Value = 3.6
$ws.Range("Q127").Value = 3
$ws.Range("Q129").Value = 3.75
$ws.Range("R127").Value = -0.25
$ws.Range("R129").Value = -0.5
$ws.Range("S127").Value = 1.85
$ws.Range("S129").Value = 1.825
$ws.Range("T127").Value = 1.95
$ws.Range("T129").Value = 1.975
$ws.Range("U127").Value = 3.5
$ws.Range("U129").Value = 3.25
$ws.Range("V127").Value = 1.825
$ws.Range("V129").Value = 1.925
$ws.Range("W127").Value = 1.975
$ws.Range("W129").Value = 1.875
$ws.Range("X127").Value = 1.05
$ws.Range("X129").Value = -1
$ws.Range("Y127").Value = -1
$ws.Range("Y129").Value = 2.6
$ws.Range("Z127").Value = -1
$ws.Range("Z129").Value = -1
$ws.Range("AA127").Value = 0.8500000000000001
$ws.Range("AA129").Value = -1
$ws.Range("AB127").Value = -1
$ws.Range("AB129").Value = 0.9750000000000001
$ws.Range("AC127").Value = -1
$ws.Range("AC129").Value = 0.925
$ws.Range("AD127").Value = 0.9750000000000001
$ws.Range("AD129").Value = -1
